$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from template rows so the new rows reuse the existing
# cell styles (A column = date style; B:Q = style "1" for rows 217-221,
# style "11" for rows 222-227), matching the source workbook pattern.
$ws.Range("A213:Q213").Copy()
$ws.Range("A217:Q221").PasteSpecial(-4122)
$ws.Range("A216:Q216").Copy()
$ws.Range("A222:Q227").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 217
$ws.Range("A217").Value = 44167
$ws.Range("B217").Value = 1007130
$ws.Range("C217").Value = 144138
$ws.Range("D217").Value = 861388
$ws.Range("E217").Value = 1604
$ws.Range("F217").Value = 2299
$ws.Range("G217").Value = 138430
$ws.Range("H217").Value = 196
$ws.Range("I217").Value = 21
$ws.Range("J217").Value = 39
$ws.Range("K217").Value = 0
$ws.Range("L217").Value = 5
$ws.Range("M217").Value = 876
$ws.Range("N217").Value = 1219
$ws.Range("O217").Value = 39
$ws.Range("P217").Value = 123
$ws.Range("Q217").Value = 42

# Row 218
$ws.Range("A218").Value = 44168
$ws.Range("B218").Value = 1011938
$ws.Range("C218").Value = 144359
$ws.Range("D218").Value = 865558
$ws.Range("E218").Value = 2021
$ws.Range("F218").Value = 2303
$ws.Range("G218").Value = 138574
$ws.Range("H218").Value = 198
$ws.Range("I218").Value = 24
$ws.Range("J218").Value = 40
$ws.Range("K218").Value = 0
$ws.Range("L218").Value = 7
$ws.Range("M218").Value = 878
$ws.Range("N218").Value = 1221
$ws.Range("O218").Value = 39
$ws.Range("P218").Value = 123
$ws.Range("Q218").Value = 42

# Row 219
$ws.Range("A219").Value = 44169
$ws.Range("B219").Value = 1017193
$ws.Range("C219").Value = 144598
$ws.Range("D219").Value = 870383
$ws.Range("E219").Value = 2212
$ws.Range("F219").Value = 2304
$ws.Range("G219").Value = 138629
$ws.Range("H219").Value = 208
$ws.Range("I219").Value = 20
$ws.Range("J219").Value = 39
$ws.Range("K219").Value = 0
$ws.Range("L219").Value = 7
$ws.Range("M219").Value = 879
$ws.Range("N219").Value = 1221
$ws.Range("O219").Value = 39
$ws.Range("P219").Value = 123
$ws.Range("Q219").Value = 42

# Row 220
$ws.Range("A220").Value = 44170
$ws.Range("B220").Value = 1021704
$ws.Range("C220").Value = 144816
$ws.Range("D220").Value = 874575
$ws.Range("E220").Value = 2313
$ws.Range("F220").Value = 2304
$ws.Range("G220").Value = 138833
$ws.Range("H220").Value = 206
$ws.Range("I220").Value = 21
$ws.Range("J220").Value = 40
$ws.Range("K220").Value = 0
$ws.Range("L220").Value = 7
$ws.Range("M220").Value = 879
$ws.Range("N220").Value = 1221
$ws.Range("O220").Value = 39
$ws.Range("P220").Value = 123
$ws.Range("Q220").Value = 42

# Row 221
$ws.Range("A221").Value = 44171
$ws.Range("B221").Value = 1024792
$ws.Range("C221").Value = 144996
$ws.Range("D221").Value = 877272
$ws.Range("E221").Value = 2524
$ws.Range("F221").Value = 2305
$ws.Range("G221").Value = 139081
$ws.Range("H221").Value = 211
$ws.Range("I221").Value = 34
$ws.Range("J221").Value = 41
$ws.Range("K221").Value = 0
$ws.Range("L221").Value = 7
$ws.Range("M221").Value = 879
$ws.Range("N221").Value = 1222
$ws.Range("O221").Value = 39
$ws.Range("P221").Value = 123
$ws.Range("Q221").Value = 42

# Row 222
$ws.Range("A222").Value = 44172
$ws.Range("B222").Value = 1025580
$ws.Range("C222").Value = 145045
$ws.Range("D222").Value = 878044
$ws.Range("E222").Value = 2491
$ws.Range("F222").Value = 2305
$ws.Range("G222").Value = 139081
$ws.Range("H222").Value = 207
$ws.Range("I222").Value = 39
$ws.Range("J222").Value = 40
$ws.Range("K222").Value = 0
$ws.Range("L222").Value = 6
$ws.Range("M222").Value = 879
$ws.Range("N222").Value = 1222
$ws.Range("O222").Value = 39
$ws.Range("P222").Value = 123
$ws.Range("Q222").Value = 42

# Row 223
$ws.Range("A223").Value = 44173
$ws.Range("B223").Value = 1027730
$ws.Range("C223").Value = 145172
$ws.Range("D223").Value = 879789
$ws.Range("E223").Value = 2769
$ws.Range("F223").Value = 2308
$ws.Range("G223").Value = 139292
$ws.Range("H223").Value = 204
$ws.Range("I223").Value = 28
$ws.Range("J223").Value = 41
$ws.Range("K223").Value = 0
$ws.Range("L223").Value = 5
$ws.Range("M223").Value = 881
$ws.Range("N223").Value = 1223
$ws.Range("O223").Value = 39
$ws.Range("P223").Value = 123
$ws.Range("Q223").Value = 42

# Row 224
$ws.Range("A224").Value = 44174
$ws.Range("B224").Value = 1029030
$ws.Range("C224").Value = 145207
$ws.Range("D224").Value = 881051
$ws.Range("E224").Value = 2772
$ws.Range("F224").Value = 2309
$ws.Range("G224").Value = 139489
$ws.Range("H224").Value = 216
$ws.Range("I224").Value = 28
$ws.Range("J224").Value = 43
$ws.Range("K224").Value = 0
$ws.Range("L224").Value = 6
$ws.Range("M224").Value = 881
$ws.Range("N224").Value = 1224
$ws.Range("O224").Value = 39
$ws.Range("P224").Value = 123
$ws.Range("Q224").Value = 42

# Row 225
$ws.Range("A225").Value = 44175
$ws.Range("B225").Value = 1033562
$ws.Range("C225").Value = 145505
$ws.Range("D225").Value = 885366
$ws.Range("E225").Value = 2691
$ws.Range("F225").Value = 2311
$ws.Range("G225").Value = 139687
$ws.Range("H225").Value = 212
$ws.Range("I225").Value = 29
$ws.Range("J225").Value = 43
$ws.Range("K225").Value = 0
$ws.Range("L225").Value = 6
$ws.Range("M225").Value = 882
$ws.Range("N225").Value = 1225
$ws.Range("O225").Value = 39
$ws.Range("P225").Value = 123
$ws.Range("Q225").Value = 42

# Row 226
$ws.Range("A226").Value = 44176
$ws.Range("B226").Value = 1037631
$ws.Range("C226").Value = 145709
$ws.Range("D226").Value = 888851
$ws.Range("E226").Value = 3071
$ws.Range("F226").Value = 2314
$ws.Range("G226").Value = 139871
$ws.Range("H226").Value = 222
$ws.Range("I226").Value = 27
$ws.Range("J226").Value = 43
$ws.Range("K226").Value = 0
$ws.Range("L226").Value = 7
$ws.Range("M226").Value = 884
$ws.Range("N226").Value = 1226
$ws.Range("O226").Value = 39
$ws.Range("P226").Value = 123
$ws.Range("Q226").Value = 42

# Row 227
$ws.Range("A227").Value = 44177
$ws.Range("B227").Value = 1042534
$ws.Range("C227").Value = 145958
$ws.Range("D227").Value = 893012
$ws.Range("E227").Value = 3564
$ws.Range("F227").Value = 2316
$ws.Range("G227").Value = 140076
$ws.Range("H227").Value = 261
$ws.Range("I227").Value = 29
$ws.Range("J227").Value = 48
$ws.Range("K227").Value = 0
$ws.Range("L227").Value = 8
$ws.Range("M227").Value = 886
$ws.Range("N227").Value = 1226
$ws.Range("O227").Value = 39
$ws.Range("P227").Value = 123
$ws.Range("Q227").Value = 42

